# Adds the new "25. 5. 2021" survey wave to both sheets (data % and pocetR sample
# sizes), refreshes a handful of rolling sample-size counts on "pocetR", updates the
# footnote date stamps, and re-selects "data" as the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "data"   - percentages
$ws2 = $wb.Worksheets.Item(2)   # "pocetR" - sample sizes

# --- Append the new date-wave column, copying the header format from the previous one ---
$ws1.Range("AC1").Copy()
$ws1.Range("AD1").PasteSpecial(-4122)
$ws1.Range("AD1").Value = "25. 5. 2021"

$ws2.Range("AB1").Copy()
$ws2.Range("AC1").PasteSpecial(-4122)
$ws2.Range("AC1").Value = "25. 5. 2021"

# --- New column values: "data" sheet (percentages), column AD, rows 2-23 ---
$ws1.Range("AD2").Value = 0.11
$ws1.Range("AD3").Value = 0.1
$ws1.Range("AD4").Value = 0.13
$ws1.Range("AD5").Value = 0.11
$ws1.Range("AD6").Value = 0.1
$ws1.Range("AD7").Value = 0.14
$ws1.Range("AD8").Value = 0.12
$ws1.Range("AD9").Value = 0.27
$ws1.Range("AD10").Value = 0.12
$ws1.Range("AD11").Value = 0.09
$ws1.Range("AD12").Value = 0.13
$ws1.Range("AD13").Value = 0.09
$ws1.Range("AD14").Value = 0.21
$ws1.Range("AD15").Value = 0.14
$ws1.Range("AD16").Value = 0.1
$ws1.Range("AD17").Value = 0.2
$ws1.Range("AD18").Value = 0.12
$ws1.Range("AD19").Value = 0.08
$ws1.Range("AD20").Value = 0.13
$ws1.Range("AD21").Value = 0.08
$ws1.Range("AD22").Value = 0.09
$ws1.Range("AD23").Value = 0.19

# --- New column values: "pocetR" sheet (sample sizes), column AC, rows 2-23 ---
$ws2.Range("AC2").Value = 1975
$ws2.Range("AC3").Value = 956
$ws2.Range("AC4").Value = 1019
$ws2.Range("AC5").Value = 269
$ws2.Range("AC6").Value = 687
$ws2.Range("AC7").Value = 313
$ws2.Range("AC8").Value = 706
$ws2.Range("AC9").Value = 168
$ws2.Range("AC10").Value = 303
$ws2.Range("AC11").Value = 376
$ws2.Range("AC12").Value = 352
$ws2.Range("AC13").Value = 776
$ws2.Range("AC14").Value = 186
$ws2.Range("AC15").Value = 421
$ws2.Range("AC16").Value = 1368
$ws2.Range("AC17").Value = 225
$ws2.Range("AC18").Value = 738
$ws2.Range("AC19").Value = 639
$ws2.Range("AC20").Value = 249
$ws2.Range("AC21").Value = 512
$ws2.Range("AC22").Value = 844
$ws2.Range("AC23").Value = 619

# --- Refreshed rolling sample-size counts in existing column AB of "pocetR" ---
$ws2.Range("AB2").Value = 2029
$ws2.Range("AB3").Value = 980
$ws2.Range("AB4").Value = 1049
$ws2.Range("AB5").Value = 271
$ws2.Range("AB8").Value = 734
$ws2.Range("AB10").Value = 311
$ws2.Range("AB11").Value = 385
$ws2.Range("AB12").Value = 364
$ws2.Range("AB15").Value = 425
$ws2.Range("AB16").Value = 1386
$ws2.Range("AB17").Value = 239
$ws2.Range("AB19").Value = 681
$ws2.Range("AB20").Value = 253
$ws2.Range("AB23").Value = 737

# --- Footnote date stamps (last row, column A) on both sheets ---
$ws1.Range("A24").Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 1. 6. 2021"
$ws2.Range("A24").Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 6. 2021"

# --- "data" becomes the active/selected sheet (was "pocetR") ---
$ws1.Activate()
